$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 5; existing rows 5-9 shift down to 7-11.
$ws.Rows.Item(5).Resize(2).Insert()

# New row 5 = copy of row 2 (Primera) but with date 2022-07-22 (44764)
$ws.Range("A5:R5").Value = $ws.Range("A2:R2").Value2
$ws.Range("D5").Value = 44764

# New row 6 = copy of row 3 (Segunda) but with date 2022-07-22 (44764)
$ws.Range("A6:R6").Value = $ws.Range("A3:R3").Value2
$ws.Range("D6").Value = 44764

# Match the date-number-format style used by the other date cells in column D.
$ws.Range("D5").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("D6").NumberFormat = $ws.Range("D3").NumberFormat
